# Fix board outline Gerbers
# Update resistor comment values in the BOM sheet:
#   "5.1k, 0603" -> "5.1kR, 0603"
#   "1.2k, 0603" -> "1.2kR, 0603"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bill of Materials-Magnetic-Sequ")

$ws.Range("C9").Value = "'5.1kR, 0603"
$ws.Range("C11").Value = "'1.2kR, 0603"
